$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12 (geometry / createreferenceframe): was TODO/easy in D12:E12 -> now documented,
# so clear the TODO columns and record the Python function name in column C.
$ws.Range("D12:E12").ClearContents()
$ws.Range("C12").Value = "create_reference_frame"

# Row 16 (geometry / getglobalcoordinates): same pattern.
$ws.Range("D16:E16").ClearContents()
$ws.Range("C16").Value = "get_global_coordinates"

# Row 17 (geometry / getlocalcoordinates): same pattern.
$ws.Range("D17:E17").ClearContents()
$ws.Range("C17").Value = "get_local_coordinates"

# Row 68 (TimeSeries / timenormalize): was TODO/moderate with a note in D68:F68 ->
# now implemented as ktk.cycles.time_normalize, recorded in column C.
$ws.Range("D68:F68").ClearContents()
$ws.Range("C68").Value = "ktk.cycles.time_normalize"

# Move the active selection to reflect where the author was working (row 67 area).
$ws.Range("A67").Select()
